# Test Plan for Professional Practices.xlsx -- add two test case rows
# (TC.001 "Opening App" and TC.002 "Markers") to the test-plan table on
# the first worksheet, and flip the pass/fail markers on rows 5 and 6.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 4: TC.001 / Opening App -------------------------------------
$ws.Cells.Item(4, 2).Value = 1
$ws.Cells.Item(4, 3).Value = "TC.001"
$ws.Cells.Item(4, 4).Value = "Opening App"
$ws.Cells.Item(4, 5).Value = "GUI"
$ws.Cells.Item(4, 6).Value = "1. Launch Home page"
$ws.Cells.Item(4, 7).Value = "Google Maps will be displayed with the users location"
$ws.Cells.Item(4, 8).Value = "Google maps displayed, accurately displaying the users location"
$ws.Cells.Item(4, 9).Value = "Google maps displayed, accurately displaying the users location"

# --- Row 5: TC.002 / Markers ------------------------------------------
$ws.Cells.Item(5, 2).Value = 2
$ws.Cells.Item(5, 3).Value = "TC.002"
$ws.Cells.Item(5, 4).Value = "Markers"
$ws.Cells.Item(5, 5).Value = "GUI"
$ws.Cells.Item(5, 6).Value = "1. Launch Home page 2.Review Markers Positions"
$ws.Cells.Item(5, 7).Value = "Users position displayed with markers of potholes or accidents close by"
$ws.Cells.Item(5, 8).Value = "The users location will be displayed in the centre of the screen, and markers in the surrounding area"
$ws.Cells.Item(5, 9).Value = "The users location will be displayed in the centre of the screen, and markers in the surrounding area"
$ws.Cells.Item(5, 10).Value = "Pass"

# --- Row 6: flip the result to a (lowercase) fail ----------------------
$ws.Cells.Item(6, 10).Value = "fail"

# --- Row heights now that rows 4/5 hold multi-line wrapped text --------
$ws.Rows.Item(4).RowHeight = 57.6
$ws.Rows.Item(5).RowHeight = 86.4

# --- Selection / scroll position, as left by the author on save --------
$ws.Application.GoTo($ws.Range("K7"), $false)
$ws.Application.ActiveWindow.ScrollRow = 1
